# Apply the "Drop in files from RMI script" edit:
#  - Remove the "Texas Data" worksheet entirely.
#  - Restore the original (pre-fix) hydrogen-from-natural-gas efficiency
#    formula on the HPEbP sheet (cell B3) to include the waste-heat term.

$wb = $excel.ActiveWorkbook

# Delete the "Texas Data" sheet
$texasSheet = $wb.Worksheets.Item("Texas Data")
$texasSheet.Delete()

# Restore the old formula in HPEbP!B3 (was =118/(162+2))
$hpebp = $wb.Worksheets.Item("HPEbP")
$hpebp.Range("B3").Formula = "=118/(162+2+46)"

# Restore the prior window/selection state (the "About" sheet was the
# active tab in this older copy of the workbook).
$ieaData = $wb.Worksheets.Item("IEA Data")
$ieaData.Range("A30").Select()

$hpebp.Range("H17:H18").Select()

$about = $wb.Worksheets.Item("About")
$about.Activate()
$about.Range("B16").Select()
